$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lat_long")
$ws.Activate()
$win = $excel.ActiveWindow
$asv = $win.ActiveSheetView
$props = $asv | Get-Member -MemberType Property
Write-Output $props
